$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 9 grows taller to accommodate the three new wrapped/rotated headers.
$ws.Rows(9).RowHeight = 179.25

# K9:M9 used to hold the plain numbers 6/7/8 (column index hints); they now
# carry the same kind of rotated, wrapped text label used by F9:J9.
$ws.Range("K9").Value = "จัดทำ Slide (Google Slide + Extension)(10)"
$ws.Range("L9").Value = "จัดทำบทเรียนด้วย Google Doc + Extension (10)"
$ws.Range("M9").Value = "ทำงานวิจัยด้วย ChatDOC (10)"

# Match the 90-degree rotated style already used by the neighboring F9:J9
# header cells (this resolves to the same cell style as those cells).
$ws.Range("K9:M9").Orientation = 90

# Update the view state: scroll position and active selection.
$ws.Range("K9").Select()
$excel.ActiveWindow.ScrollRow = 4
$excel.ActiveWindow.ScrollColumn = 1
